# Weekly Time Record - time-in/time-out updates for 30-May-2016 week
# (Tue/Wed/Thu/Fri rows) + follow-on UI state (selection, print area,
# column widths) that Excel/LibreOffice touches as a side effect of
# editing & re-saving the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Time in / time out entries -------------------------------------------------
# Row 15 (Wednesday): 6:00 PM - 7:00 PM
$ws.Range("C15").Value2 = 0.75
$ws.Range("D15").Value2 = 0.791666666666667

# Row 17 (Friday): 6:00 PM - 7:00 PM
$ws.Range("C17").Value2 = 0.75
$ws.Range("D17").Value2 = 0.791666666666667

# Row 18 (Saturday): 11:00 AM - 12:00 PM
$ws.Range("C18").Value2 = 0.458333333333333
$ws.Range("D18").Value2 = 0.5

# Row 19 (Sunday): 2:00 PM - 9:00 PM
$ws.Range("C19").Value2 = 0.583333333333333
$ws.Range("D19").Value2 = 0.875

# --- Re-assert the print area (source file accumulates a fresh
#     _xlnm.Print_Area_* defined name each time this is done) -------------------
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# --- Minor column-width settle from the edit/save round-trip -------------------
$ws.Columns("B").ColumnWidth = 11.86
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 9.56
$ws.Columns("H").ColumnWidth = 11.04
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 9.56
$ws.StandardWidth = 6.73

# --- Active cell moved to C18 after entering the data ---------------------------
[void]$ws.Range("C18").Select()
